$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I, shifting the existing average columns
# (I:L -> J:M) to the right.
$ws.Columns("I:I").Insert()

# New column I header mirrors the "background" label used in column B,
# matching the pattern of the other average-header cells (I1:L1 originally).
# Copy J1 first so I1 inherits the bold/border/centered header style, then
# copy B1's value (and string identity) on top without touching formatting.
$ws.Range("J1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("I1").PasteSpecial(-4163)

# New column I row 2 formula: average of the "background" column (B).
$ws.Range("I2").Formula = "=AVERAGE(B2:B162)"

# Leave the shared-formula group for J2:M2 as individual (non-shared)
# formulas referencing their original source columns. Clear first so the
# engine does not keep treating them as one shared-formula group.
$ws.Range("J2:M2").ClearContents()
$ws.Range("J2").Formula = "=AVERAGE(D2:D162)"
$ws.Range("K2").Formula = "=AVERAGE(E2:E162)"
$ws.Range("L2").Formula = "=AVERAGE(F2:F162)"
$ws.Range("M2").Formula = "=AVERAGE(G2:G162)"

# Reflect the final cell selection recorded in the saved workbook.
$ws.Range("G5").Select()

$wb.Save()
